# ---------------------------------------------------------------------------
# Reproduces the commit "added extent repoting and data provider":
#   * removes the "tc_01" test-case sheet
#   * bumps the password value used by tc_02 / tc_03
#   * refreshes the "testData" sheet's tc_01-labelled block to describe tc_03
#     instead, and renames its last block to the new test case name
#   * appends a brand-new "validateFBLoginWIthValidCred" sheet with its own
#     username/password rows (one of them hyperlinked, like the existing
#     testData mail-to links)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# 1) Drop the tc_01 sheet entirely -- every other tab/file shifts down one
#    slot automatically when this is saved.
$wb.Worksheets.Item("tc_01").Delete() | Out-Null

# 2) tc_02 : password value changes, and the header row gets formatted as
#    text (mirrors the "Format Cells -> Text" applied by the original author
#    before they retyped the password).
$tc02 = $wb.Worksheets.Item("tc_02")
$tc02.Range("A1:B1").NumberFormat = "@"
$tc02.Range("B2").Value = 4545645
$tc02.Rows("1:1").Select() | Out-Null

# 3) tc_03 : same treatment, different password value.
$tc03 = $wb.Worksheets.Item("tc_03")
$tc03.Range("A1:B1").NumberFormat = "@"
$tc03.Range("B2").Value = 56767567
$tc03.Rows("1:1").Select() | Out-Null

# 4) tc_04 / tc_05 / tc_06 keep their own data untouched; only their backing
#    part numbers shift because of the tc_01 deletion above.

# 5) testData : the tc_01 example block becomes a tc_03 example block (only
#    two columns now, email/firstname/lastname go away), and the final
#    block is renamed for the newly added test case.
$testData = $wb.Worksheets.Item("testData")
$testData.Range("A5").Value = "tc_03"
$testData.Range("C6:E6").ClearContents() | Out-Null
$testData.Range("A11").Value = "validateFBLoginWIthValidCred"
$testData.Range("A13").Select() | Out-Null

# 6) Add the new sheet at the very end of the workbook for the new test
#    case, with its own username/password sample rows.
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $last)
$newSheet.Name = "validateFBLoginWIthValidCred"

$newSheet.Range("A1").Value = "username"
$newSheet.Range("B1").Value = "password"
$newSheet.Range("A1:B1").NumberFormat = "@"

$newSheet.Range("A2").Value = "rahul.jhajava@kljflkjds"
$newSheet.Range("B2").Value = "kjkasldj"
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Hyperlinks.Add($newSheet.Range("A2"), "mailto:rahul.jhajava@kljflkjds") | Out-Null
$newSheet.Range("A2").NumberFormat = "@"

$newSheet.Range("A3").Value = "sdsdf"
$newSheet.Range("B3").Value = "asdsad"

$newSheet.Columns.Item(1).ColumnWidth = 21.7109375
$newSheet.Columns.Item(2).ColumnWidth = 9.42578125

$newSheet.PageSetup.Orientation = 1

# 7) Leave testData as the active tab, matching the saved workbook view.
$testData.Activate() | Out-Null
